# Localisation-MED scenario update:
#  - "TFM_FILL" flag renamed to "~TFM_FILL" (disable the TFM_FILL row so VEDA
#    skips it, per the convention of prefixing tags with "~" to comment them out)
#  - Electricity capital-cost (NCAP_COST) data refreshed: a new 2017 base-year
#    value was inserted ahead of the old series (old 2017 -> new 2020 slot),
#    and the 2025/2030/2040/2050 projections were updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Flag rename -----------------------------------------------------------
$ws.Range("U28").Value = "~TFM_FILL"

# --- Row 37: Solar Central Receiver 09 hrs storage (ERSOLTC09-N) -----------
$ws.Range("X37").Value = 81278.984337530506
$ws.Range("AC37").Value = 81279.984337530506
$ws.Range("AH37").Value = 63100.036393056798
$ws.Range("AM37").Value = 52213.379980116501
$ws.Range("AS37").Value = 53456.5556939288
$ws.Range("AX37").Value = 54699.731407741099

# --- Row 38: Solar PV Fixed (ERSOLPCF-N) ------------------------------------
$ws.Range("X38").Value = 48762.532743362797
$ws.Range("AC38").Value = 19505.013097345101
$ws.Range("AH38").Value = 13222.764078584099
$ws.Range("AM38").Value = 12787.684470796399
$ws.Range("AS38").Value = 11630.3803637168
$ws.Range("AX38").Value = 10513.896293805299

# --- Row 39: Solar PV tracking (ERSOLPCT-N) ---------------------------------
$ws.Range("X39").Value = 51982.699999999997
$ws.Range("AC39").Value = 20793.080000000002
$ws.Range("AH39").Value = 14095.965480000001
$ws.Range("AM39").Value = 13632.154200000001
$ws.Range("AS39").Value = 12398.424349999999
$ws.Range("AX39").Value = 11208.2102

# --- Row 40: Wind (ERWNDH-N) -------------------------------------------------
$ws.Range("X40").Value = 32295.527999999998
$ws.Range("AC40").Value = 21530.351999999999
$ws.Range("AH40").Value = 17770.180919999999
$ws.Range("AM40").Value = 18133.392899999999
$ws.Range("AS40").Value = 18241.690050000001
$ws.Range("AX40").Value = 18341.129400000002

# --- Row 48: Biomass municipal waste (ERBIO-N) ------------------------------
$ws.Range("AC48").Value = 38195.713600000003
$ws.Range("AH48").Value = 24657.861225652799
$ws.Range("AM48").Value = 21276.840861764798
$ws.Range("AS48").Value = 18424.264945496401
$ws.Range("AX48").Value = 15417.019777003599

# --- Row 49: Solar PV rooftop Agriculture (ERSOLPRA-N) ----------------------
$ws.Range("X49").Value = 56564.537982300899
$ws.Range("AC49").Value = 22625.8151929204
$ws.Range("AH49").Value = 15338.406331157499
$ws.Range("AM49").Value = 14833.713986123899
$ws.Range("AS49").Value = 13491.241221911499
$ws.Range("AX49").Value = 12196.119700814201

# --- Row 50: Solar PV rooftop Mining (ERSOLPRM-N) ---------------------------
$ws.Range("X50").Value = 56564.537982300899
$ws.Range("AC50").Value = 22625.8151929204
$ws.Range("AH50").Value = 15338.406331157499
$ws.Range("AM50").Value = 14833.713986123899
$ws.Range("AS50").Value = 13491.241221911499
$ws.Range("AX50").Value = 12196.119700814201

# --- Row 51: Solar PV rooftop commercial (ERSOLPRC-N) -----------------------
$ws.Range("X51").Value = 56564.537982300899
$ws.Range("AC51").Value = 22625.8151929204
$ws.Range("AH51").Value = 15338.406331157499
$ws.Range("AM51").Value = 14833.713986123899
$ws.Range("AS51").Value = 13491.241221911499
$ws.Range("AX51").Value = 12196.119700814201

# --- Row 52: Solar PV rooftop residential (ERSOLPRR-N) ----------------------
$ws.Range("X52").Value = 85334.432300885004
$ws.Range("AC52").Value = 34133.772920354
$ws.Range("AH52").Value = 23139.837137522201
$ws.Range("AM52").Value = 22378.447823893799
$ws.Range("AS52").Value = 20353.1656365044
$ws.Range("AX52").Value = 18399.318514159299

# --- Row 53: Solar PV rooftop Industry (ERSOLPRI-N) -------------------------
$ws.Range("X53").Value = 56564.537982300899
$ws.Range("AC53").Value = 22625.8151929204
$ws.Range("AH53").Value = 15338.406331157499
$ws.Range("AM53").Value = 14833.713986123899
$ws.Range("AS53").Value = 13491.241221911499
$ws.Range("AX53").Value = 12196.119700814201

# --- Row 54: Transmission (ETRANS) ------------------------------------------
$ws.Range("X54").Value = 1637.60737116

# --- Row 55: dummy tech tracking investment costs (ETRANSDUM) --------------
$ws.Range("X55").Value = 337.44

# --- Row 56: Agricultural Electricity (XAGRELC) -----------------------------
$ws.Range("X56").Value = 19046.524215158599

# --- Row 57: Commercial Electricity (XCOMELC) -------------------------------
$ws.Range("X57").Value = 9523.2621384296508

# --- Row 58: Industrial CP-Electricity (XICPELC) ----------------------------
$ws.Range("X58").Value = 4526.8304959535499

# --- Row 59: Industrial FA-Electricity (XIFAELC) ----------------------------
$ws.Range("X59").Value = 4526.8304959535499

# --- Row 60: Industrial FB-Electricity (XIFBELC) ----------------------------
$ws.Range("X60").Value = 7272.9396030149201

# --- Row 61: Industrial IS-Electricity (XIISELC) ----------------------------
$ws.Range("X61").Value = 4526.8304959535499

# --- Row 62: Industrial MI-Electricity (XIMIELC) ----------------------------
$ws.Range("X62").Value = 4526.8304959535499

# --- Row 63: Industrial NF-Electricity (XINFELC) ----------------------------
$ws.Range("X63").Value = 4526.8304959535499

# --- Row 64: Industrial NM-Electricity (XINMELC) ----------------------------
$ws.Range("X64").Value = 4526.8304959535499

# --- Row 65: Industrial OT-Electricity (XIOTELC) ----------------------------
$ws.Range("X65").Value = 7272.9396030149201

# --- Row 66: Industrial PP-Electricity (XIPPELC) ----------------------------
$ws.Range("X66").Value = 7272.9396030149201

# --- Row 67: Residential Electricity (XRESELC) ------------------------------
$ws.Range("X67").Value = 19046.524215158599

# --- Row 68: Transport Electricity (XTRAELC) --------------------------------
$ws.Range("X68").Value = 9523.2621384296508

# --- Row 69: Supply Sector Electricity (XUPSELC) ----------------------------
$ws.Range("X69").Value = 4761.6310692148199

# --- Restore the cursor position that was active when the file was saved ---
$ws.Range("U29").Select()
